$d = $word.ActiveDocument

# Locate the paragraph that contains the "GIS & Geospatial Analysis Consulting"
# heading line that immediately precedes the Siege Analytics bullet list.
$idx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r", "`n", "`a")
    if ($t -eq "GIS & Geospatial Analysis Consulting") {
        $idx = $i
        break
    }
}

if ($idx -eq -1) {
    throw "Could not find the 'GIS & Geospatial Analysis Consulting' paragraph"
}

$bullets = @(
    [char]0x2022 + " Utilized ESRI Arc Suite and OSGeo technology to map and analyze 50,000+ electoral boundaries across federal, state, and local levels",
    [char]0x2022 + " Applied geospatial analysis to uncover demographic miscoding affecting 2,000+ precincts nationwide",
    [char]0x2022 + " Developed boundary estimation tools enabling smaller organizations to conduct sophisticated redistricting analysis"
)

$insertAfter = $idx
foreach ($bullet in $bullets) {
    $anchor = $d.Paragraphs.Item($insertAfter)
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertAfter + 1)
    $newPara.Range.Text = $bullet
    $insertAfter = $insertAfter + 1
}
